$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 9
$ws.Range("B6").Value = 11

$ws.Range("D8").Select()
